$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "dSF" column (F) values for the rows that changed.
$ws.Range("F5").Value = -9
$ws.Range("F8").Value = -11
$ws.Range("F12").Value = -4
$ws.Range("F15").Value = -3
$ws.Range("F16").Value = -7
